# Auto-generated Excel COM-interop script to apply numeric corrections
# to the LeveProfit calculation sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 435.6316
$ws.Range("I92").Value = 435.6316
$ws.Range("K92").Value = 435.6316
$ws.Range("M92").Value = 812.3684000000001
# Row 113
$ws.Range("H113").Value = 113367.22
$ws.Range("J113").Value = 2200
$ws.Range("L113").Value = 2200
$ws.Range("N113").Value = -8708
# Row 137
$ws.Range("H137").Value = 1214.4375
$ws.Range("I137").Value = 1245.6154
$ws.Range("J137").Value = 1079.3334
$ws.Range("K137").Value = 3736.8462
$ws.Range("L137").Value = 3238.0002
$ws.Range("M137").Value = -1186.8462
$ws.Range("N137").Value = -8338.0002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 39655.383
$ws.Range("I2").Value = 1109.25
$ws.Range("J2").Value = 101329.2
$ws.Range("K2").Value = 1109.25
$ws.Range("L2").Value = 101329.2
$ws.Range("M2").Value = -996.25
$ws.Range("N2").Value = -101555.2
# Row 45
$ws.Range("H45").Value = 1633.0834
$ws.Range("I45").Value = 1762.6666
$ws.Range("J45").Value = 1503.5
$ws.Range("K45").Value = 1762.6666
$ws.Range("L45").Value = 1503.5
$ws.Range("M45").Value = -1385.6666
$ws.Range("N45").Value = -2257.5
# Row 61
$ws.Range("H61").Value = 1837.875
$ws.Range("I61").Value = 1200.5
$ws.Range("J61").Value = 3750
$ws.Range("K61").Value = 1200.5
$ws.Range("L61").Value = 3750
$ws.Range("M61").Value = -988.5
$ws.Range("N61").Value = -4174
# Row 116
$ws.Range("H116").Value = 39655.383
$ws.Range("I116").Value = 1109.25
$ws.Range("J116").Value = 101329.2
$ws.Range("K116").Value = 1109.25
$ws.Range("L116").Value = 101329.2
$ws.Range("M116").Value = 1184.75
$ws.Range("N116").Value = -105917.2
# Row 136
$ws.Range("H136").Value = 1837.875
$ws.Range("I136").Value = 1200.5
$ws.Range("J136").Value = 3750
$ws.Range("K136").Value = 3601.5
$ws.Range("L136").Value = 11250
$ws.Range("M136").Value = -1051.5
$ws.Range("N136").Value = -16350

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 39655.383
$ws.Range("I3").Value = 1109.25
$ws.Range("J3").Value = 101329.2
$ws.Range("K3").Value = 1109.25
$ws.Range("L3").Value = 101329.2
$ws.Range("M3").Value = -995.25
$ws.Range("N3").Value = -101557.2
# Row 134
$ws.Range("H134").Value = 2207.3823
$ws.Range("I134").Value = 2435.2593
$ws.Range("J134").Value = 1328.4286
$ws.Range("K134").Value = 7305.777900000001
$ws.Range("L134").Value = 3985.2858
$ws.Range("M134").Value = -4770.777900000001
$ws.Range("N134").Value = -9055.2858

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 20597.055
$ws.Range("I31").Value = 25690.805
$ws.Range("J31").Value = 4532.154
$ws.Range("K31").Value = 25690.805
$ws.Range("L31").Value = 4532.154
$ws.Range("M31").Value = -25395.805
$ws.Range("N31").Value = -5122.154
# Row 34
$ws.Range("H34").Value = 20597.055
$ws.Range("I34").Value = 25690.805
$ws.Range("J34").Value = 4532.154
$ws.Range("K34").Value = 25690.805
$ws.Range("L34").Value = 4532.154
$ws.Range("M34").Value = -25488.805
$ws.Range("N34").Value = -4936.154
# Row 58
$ws.Range("H58").Value = 7594.795
$ws.Range("I58").Value = 860.3200000000001
$ws.Range("J58").Value = 19620.643
$ws.Range("K58").Value = 860.3200000000001
$ws.Range("L58").Value = 19620.643
$ws.Range("M58").Value = -657.3200000000001
$ws.Range("N58").Value = -20026.643
# Row 132
$ws.Range("H132").Value = 60004240
$ws.Range("I132").Value = 62504460
$ws.Range("J132").Value = 55559404
$ws.Range("K132").Value = 187513380
$ws.Range("L132").Value = 166678212
$ws.Range("M132").Value = -187510850
$ws.Range("N132").Value = -166683272
# Row 134
$ws.Range("H134").Value = 1550.9286
$ws.Range("I134").Value = 1643.5834
$ws.Range("J134").Value = 995
$ws.Range("K134").Value = 4930.7502
$ws.Range("L134").Value = 2985
$ws.Range("M134").Value = -2395.7502
$ws.Range("N134").Value = -8055
# Row 136
$ws.Range("H136").Value = 7594.795
$ws.Range("I136").Value = 860.3200000000001
$ws.Range("J136").Value = 19620.643
$ws.Range("K136").Value = 2580.96
$ws.Range("L136").Value = 58861.929
$ws.Range("M136").Value = -30.96000000000004
$ws.Range("N136").Value = -63961.929

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 816.53
$ws.Range("I131").Value = 639.5714
$ws.Range("J131").Value = 829.8495
$ws.Range("K131").Value = 1918.7142
$ws.Range("L131").Value = 2489.5485
$ws.Range("M131").Value = 3121.2858
$ws.Range("N131").Value = -12569.5485

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 52
$ws.Range("H52").Value = 12749.75
$ws.Range("J52").Value = 12749.75
$ws.Range("L52").Value = 12749.75
$ws.Range("N52").Value = -13267.75
# Row 96
$ws.Range("H96").Value = 14650
$ws.Range("J96").Value = 14650
$ws.Range("L96").Value = 14650
$ws.Range("N96").Value = -20142
# Row 113
$ws.Range("H113").Value = 1699.5454
$ws.Range("I113").Value = 1424.75
$ws.Range("J113").Value = 1856.5714
$ws.Range("K113").Value = 1424.75
$ws.Range("L113").Value = 1856.5714
$ws.Range("M113").Value = 745.25
$ws.Range("N113").Value = -6196.5714
# Row 132
$ws.Range("H132").Value = 2368.55
$ws.Range("I132").Value = 2051.2942
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 6153.882599999999
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = -3623.882599999999
$ws.Range("N132").Value = -17559.0005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2130.261
$ws.Range("I7").Value = 1732.25
$ws.Range("J7").Value = 3040
$ws.Range("K7").Value = 1732.25
$ws.Range("L7").Value = 3040
$ws.Range("M7").Value = -1620.25
$ws.Range("N7").Value = -3264
# Row 126
$ws.Range("H126").Value = 2130.261
$ws.Range("I126").Value = 1732.25
$ws.Range("J126").Value = 3040
$ws.Range("K126").Value = 5196.75
$ws.Range("L126").Value = 9120
$ws.Range("M126").Value = -2726.75
$ws.Range("N126").Value = -14060
# Row 132
$ws.Range("H132").Value = 8611
$ws.Range("I132").Value = 10652
$ws.Range("J132").Value = 5549.5
$ws.Range("K132").Value = 31956
$ws.Range("L132").Value = 16648.5
$ws.Range("M132").Value = -29426
$ws.Range("N132").Value = -21708.5
# Row 136
$ws.Range("H136").Value = 2009.9
$ws.Range("I136").Value = 1899.8572
$ws.Range("J136").Value = 2266.6667
$ws.Range("K136").Value = 5699.571599999999
$ws.Range("L136").Value = 6800.000100000001
$ws.Range("M136").Value = -3149.571599999999
$ws.Range("N136").Value = -11900.0001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("N82").Value = 0
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("N85").Value = 0
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0
# Row 106
$ws.Range("H106").Value = 29933.334
$ws.Range("J106").Value = 29933.334
$ws.Range("L106").Value = 29933.334
$ws.Range("N106").Value = -32457.334
# Row 113
$ws.Range("H113").Value = 465.13635
$ws.Range("I113").Value = 353.375
$ws.Range("K113").Value = 1060.125
$ws.Range("M113").Value = 1109.875
# Row 136
$ws.Range("H136").Value = 1539.7941
$ws.Range("I136").Value = 661.05884
$ws.Range("J136").Value = 2418.5293
$ws.Range("K136").Value = 1983.17652
$ws.Range("L136").Value = 7255.5879
$ws.Range("M136").Value = 566.82348
$ws.Range("N136").Value = -12355.5879

